# The Chelyabinsk light-curve sheet had an empty "gap" row at row 181
# (row numbers jumped from 180 straight to 182), which was throwing off
# the integrated-intensity calculations downstream. Remove that blank
# row so the data closes up and renumbers contiguously (old row 182
# becomes the new row 181, old 183 becomes 182, ... old 341 becomes 340).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the blank row 181 - everything below shifts up by one.
$ws.Rows.Item(181).Delete()

# Touch the bottom of the sheet (this mirrors what the original
# authoring app left behind: the used range/dimension stretching all
# the way to the last sheet row with an otherwise-empty final row).
$ws.Cells.Item(1048576, 2).Style = "Normal"
$ws.Rows.Item(1048576).RowHeight = 12.8

# Restore/update the view: scrolled & selected further down the sheet.
$win = $excel.ActiveWindow
$win.ScrollRow = 163
$win.ScrollColumn = 1
$ws.Range("G187").Select()
